# Fix locator for scenario 1
# Updates the Sales Order number (CB101 -> DB104), related dates, and
# status values on both the "Order Info" and "Plan" sheets.

$wb = $excel.ActiveWorkbook
$wsOrder = $wb.Worksheets.Item("Order Info")
$wsPlan  = $wb.Worksheets.Item("Plan")

# ---- "Order Info" sheet ----

# Sales Order No.
$wsOrder.Range("C14").Value = "sDB104-2311001"

# Order Date / Target Date
$wsOrder.Range("C18").Value = 45250
$wsOrder.Range("C19").Value = "27 Nov 2023 - 03 Dec 2023"

# Forecast week labels
$wsOrder.Range("N22").Value = "04 Dec ~ 10 Dec"
$wsOrder.Range("P23").Value = 45252

# Grid rows 24-29: Purchase Order No. (column E)
$wsOrder.Range("E24").Value = "cDB104-2311001"
$wsOrder.Range("E25").Value = "cDB104-2311001"
$wsOrder.Range("E26").Value = "cDB104-2311001"
$wsOrder.Range("E27").Value = "cDB104-2311001"
$wsOrder.Range("E28").Value = "cDB104-2311001"
$wsOrder.Range("E29").Value = "cDB104-2311001"

# Status column (O): rows 24/25 Completed -> Processing
$wsOrder.Range("O24").Value = "Processing"
$wsOrder.Range("O25").Value = "Processing"

# Row 28: status moves from Processing to Completed's old slot;
# Delivered / Receiver Inbounded qty reset to 0
$wsOrder.Range("O28").Value = "Processing"
$wsOrder.Range("P28").Value = 0
$wsOrder.Range("R28").Value = 0

# ---- "Plan" sheet ----

$wsPlan.Range("B4").Value = "Sales Order No.:sDB104-2311001"

# Header date row
$wsPlan.Range("K7").Value = 45252
$wsPlan.Range("L7").Value = 45301
$wsPlan.Range("M7").Value = 45312
$wsPlan.Range("N7").Value = 45343
$wsPlan.Range("O7").Value = 45264
$wsPlan.Range("P7").Value = 45266

# Row 8 / 9: status Completed -> Processing, inbound qty moved from K to P
$wsPlan.Range("J8").Value = "Processing"
$wsPlan.Range("K8").Value = 0
$wsPlan.Range("P8").Value = 1620

$wsPlan.Range("J9").Value = "Processing"
$wsPlan.Range("K9").Value = 0
$wsPlan.Range("P9").Value = 1620

# Row 12: inbound qty moved from K to P
$wsPlan.Range("K12").Value = 0
$wsPlan.Range("P12").Value = 800
